# Update "want to go" counts (column F) for a handful of events on the
# "展览" (Exhibition) sheet and the "全部类型" (All Types) aggregate sheet.
# These correspond to a refreshed data pull (gh-pages output regenerated).

$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 (Exhibition) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F8").Value  = 14187   # was 14184
$ws1.Range("F9").Value  = 138     # was 137
$ws1.Range("F11").Value = 5697    # was 5695
$ws1.Range("F13").Value = 63      # was 62
$ws1.Range("F23").Value = 10479   # was 10478

# --- Sheet: 全部类型 (All Types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F9").Value  = 14187   # was 14184
$ws4.Range("F10").Value = 138     # was 137
$ws4.Range("F12").Value = 5697    # was 5695
$ws4.Range("F14").Value = 63      # was 62
$ws4.Range("F25").Value = 10479   # was 10478
